$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = '2026-02-09 04:48:17'
$ws.Range("O2").Value = '-4.0 °C'
$ws.Range("E3").Value = '2026-02-09 04:48:20'
$ws.Range("M3").Value = '-5.4 °C 4:28 TU'
$ws.Range("O3").Value = '-6.3 °C'
$ws.Range("E4").Value = '2026-02-09 04:48:22'
$ws.Range("O4").Value = '4.6 °C'
$ws.Range("E5").Value = '2026-02-09 04:48:24'
$ws.Range("O5").Value = '-5.5 °C'
$ws.Range("E6").Value = '2026-02-09 04:48:26'
$ws.Range("O6").Value = '6.7 °C'
$ws.Range("E7").Value = '2026-02-09 04:48:28'
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = '72%'
$ws.Range("E8").Value = '2026-02-09 04:48:31'
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '80%'
$ws.Range("O8").Value = '7.2 °C'
$ws.Range("E9").Value = '2026-02-09 04:48:33'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '85%'
$ws.Range("N9").Value = '3.8 °C 4:22 TU'
$ws.Range("O9").Value = '7.3 °C'
$ws.Range("E10").Value = '2026-02-09 04:48:36'
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '94%'
$ws.Range("N10").Value = '3.3 °C 4:29 TU'
$ws.Range("O10").Value = '6.0 °C'
$ws.Range("E11").Value = '2026-02-09 04:48:38'
$ws.Range("O11").Value = '2.1 °C'
$ws.Range("E12").Value = '2026-02-09 04:48:40'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '89%'
$ws.Range("N12").Value = '5.2 °C 4:17 TU'
$ws.Range("O12").Value = '7.7 °C'
$ws.Range("E13").Value = '2026-02-09 04:48:42'
$ws.Range("N13").Value = '-3.7 °C 4:05 TU'
$ws.Range("O13").Value = '-1.6 °C'
$ws.Range("E14").Value = '2026-02-09 04:48:45'
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = '96%'
$ws.Range("L14").Value = '18.0 km/h - 310º 4:23 TU'
$ws.Range("O14").Value = '7.5 °C'
$ws.Range("E15").Value = '2026-02-09 04:48:47'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '86%'
$ws.Range("N15").Value = '3.1 °C 4:28 TU'
$ws.Range("O15").Value = '5.9 °C'
$ws.Range("E16").Value = '2026-02-09 04:48:50'
$ws.Range("N16").Value = '-6.2 °C 4:11 TU'
$ws.Range("O16").Value = '-5.2 °C'
$ws.Range("E17").Value = '2026-02-09 04:48:52'
$ws.Range("O17").Value = '0.0 °C'
$ws.Range("E18").Value = '2026-02-09 04:48:55'
$ws.Range("I18").Value = '0.1 mm'
$ws.Range("N18").Value = '4.3 °C 4:26 TU'
$ws.Range("O18").Value = '6.6 °C'
$ws.Range("E19").Value = '2026-02-09 04:48:57'
$ws.Range("N19").Value = '2.9 °C 4:03 TU'
$ws.Range("E20").Value = '2026-02-09 04:48:59'
$ws.Range("E21").Value = '2026-02-09 04:49:02'
$ws.Range("N21").Value = '-0.5 °C 4:24 TU'
$ws.Range("O21").Value = '0.7 °C'
$ws.Range("E22").Value = '2026-02-09 04:49:04'
$ws.Range("E23").Value = '2026-02-09 04:49:06'
$ws.Range("E24").Value = '2026-02-09 04:49:09'
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = '86%'
$ws.Range("O24").Value = '4.4 °C'
$ws.Range("E25").Value = '2026-02-09 04:49:11'
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = '77%'
$ws.Range("K25").Value = '-0.1 MJ/m2'
$ws.Range("L25").Value = '17.3 km/h - 249º 4:02 TU'
$ws.Range("M25").Value = '-3.0 °C 4:05 TU'
$ws.Range("O25").Value = '-4.4 °C'
$ws.Range("E26").Value = '2026-02-09 04:49:13'
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = '91%'
$ws.Range("J26").Value = '1008.7 hPa'
$ws.Range("L26").Value = '18.0 km/h - 8º 4:03 TU'
$ws.Range("E27").Value = '2026-02-09 04:49:16'
$ws.Range("E28").Value = '2026-02-09 04:49:18'
$ws.Range("N28").Value = '1.5 °C 4:12 TU'
$ws.Range("O28").Value = '3.7 °C'
$ws.Range("E29").Value = '2026-02-09 04:49:20'
$ws.Range("N29").Value = '3.2 °C 4:22 TU'
$ws.Range("O29").Value = '5.7 °C'
$ws.Range("E30").Value = '2026-02-09 04:49:23'
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '94%'
$ws.Range("N30").Value = '4.9 °C 4:27 TU'
$ws.Range("O30").Value = '6.8 °C'
$ws.Range("E31").Value = '2026-02-09 04:49:25'
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '72%'
$ws.Range("O31").Value = '8.9 °C'
$ws.Range("E32").Value = '2026-02-09 04:49:27'
$ws.Range("E33").Value = '2026-02-09 04:49:29'
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '94%'
$ws.Range("N33").Value = '-1.5 °C 4:29 TU'
$ws.Range("O33").Value = '-0.4 °C'
$ws.Range("E34").Value = '2026-02-09 04:49:31'
$ws.Range("O34").Value = '-3.1 °C'
$ws.Range("E35").Value = '2026-02-09 04:49:34'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '66%'
$ws.Range("J35").Value = '1010.0 hPa'
$ws.Range("E36").Value = '2026-02-09 04:49:36'
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = '81%'
$ws.Range("O36").Value = '8.6 °C'
$ws.Range("E37").Value = '2026-02-09 04:49:39'
$ws.Range("J37").Value = '1009.3 hPa'
$ws.Range("N37").Value = '1.6 °C 4:29 TU'
$ws.Range("O37").Value = '3.6 °C'
$ws.Range("E38").Value = '2026-02-09 04:49:41'
$ws.Range("N38").Value = '4.8 °C 4:25 TU'
$ws.Range("O38").Value = '6.3 °C'
$ws.Range("E39").Value = '2026-02-09 04:49:43'
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = '84%'
$ws.Range("E40").Value = '2026-02-09 04:49:46'
$ws.Range("O40").Value = '-0.3 °C'
$ws.Range("E41").Value = '2026-02-09 04:49:48'
$ws.Range("E42").Value = '2026-02-09 04:49:51'
$ws.Range("N42").Value = '4.5 °C 4:12 TU'
$ws.Range("O42").Value = '6.6 °C'
$ws.Range("E43").Value = '2026-02-09 04:49:53'
$ws.Range("N43").Value = '5.9 °C 4:08 TU'
$ws.Range("E44").Value = '2026-02-09 04:49:55'
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = '91%'
$ws.Range("O44").Value = '-7.6 °C'
$ws.Range("E45").Value = '2026-02-09 04:49:57'
$ws.Range("J45").Value = '1009.9 hPa'
$ws.Range("N45").Value = '-1.4 °C 4:02 TU'
$ws.Range("O45").Value = '0.0 °C'
$ws.Range("E46").Value = '2026-02-09 04:50:00'
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '81%'
